$d = $word.ActiveDocument

# The original "FÉRIAS" run carries a stray self-closing <w:rPr/> (explicit
# run properties with nothing in them). Re-stamp that run's text via raw
# WordprocessingML so the rewritten run comes out clean (no empty <w:rPr/>),
# matching what Word's writer emits for a run with default formatting --
# while keeping the paragraph/run identity attributes (rsid, paraId, textId)
# exactly as they already are.
$feriasRange = $d.Range(0, 6)
$feriasXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r w:rsidR="6A96F447">
              <w:t>FÉRIAS</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
[void]$feriasRange.InsertXML($feriasXml)

# Position an empty (collapsed) range right before the document's final
# paragraph mark (i.e. at the very end of the "FÉRIAS" paragraph / story),
# so the new paragraph we inject becomes the new last paragraph, ahead of
# the sectPr.
$insertAt = $d.Content.End - 1
$r = $d.Range($insertAt, $insertAt)

# Build the new paragraph as a raw WordprocessingML fragment (wrapped in the
# pkg:package envelope Word's Range.InsertXML expects) so we get exact
# control over the run layout (spell-check proofErr markers + bookmark)
# that a real "type it out, let autocorrect/spellcheck tag it" edit would
# produce: "Check list guia para entrevista." in red, with "Check" and
# "list" individually bracketed by proofErr spellStart/spellEnd (Word
# flags them since they're English words in a pt-BR document), plus the
# automatic "_GoBack" bookmark Word drops at the last edit position.
$newParagraphXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>Check</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t>list</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:color w:val="FF0000"/>
              </w:rPr>
              <w:t xml:space="preserve"> guia para entrevista.</w:t>
            </w:r>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$r.InsertXML($newParagraphXml)
